function Set-TextValue($ws, $addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.136.44'
$ws.Range('E2').Value = '  -1.79%  '
$ws.Range('D3').Value = '3.113.09'
$ws.Range('E3').Value = '  -2.79%  '
$ws.Range('E4').Value = '  +0.00%  '
Set-TextValue $ws 'D5' '594.44'
$ws.Range('E5').Value = '  -0.48%  '
Set-TextValue $ws 'D6' '157.61'
$ws.Range('E6').Value = '  +2.03%  '
$ws.Range('E7').Value = '  -0.01%  '
Set-TextValue $ws 'D8' '0.543'
$ws.Range('E8').Value = '  +0.08%  '
$ws.Range('D9').Value = '3.111.45'
$ws.Range('E9').Value = '  -2.80%  '
$ws.Range('E10').Value = '  -5.18%  '
Set-TextValue $ws 'D11' '5.90'
$ws.Range('E11').Value = '  -3.43%  '
Set-TextValue $ws 'D12' '0.452'
$ws.Range('E12').Value = '  -4.04%  '
Set-TextValue $ws 'D13' '37.18'
$ws.Range('E13').Value = '  -5.75%  '
$ws.Range('E14').Value = '  -5.99%  '
$ws.Range('D15').Value = '3.626.68'
$ws.Range('E15').Value = '  -2.87%  '
$ws.Range('E16').Value = '  -1.58%  '
Set-TextValue $ws 'D17' '7.24'
$ws.Range('E17').Value = '  -2.64%  '
$ws.Range('D18').Value = '64.106.05'
$ws.Range('E18').Value = '  -1.45%  '
$ws.Range('D19').Value = '3.114.89'
$ws.Range('E19').Value = '  -2.79%  '
Set-TextValue $ws 'D20' '478.76'
$ws.Range('E20').Value = '  -1.08%  '
Set-TextValue $ws 'D21' '14.49'
$ws.Range('E21').Value = '  -4.30%  '
Set-TextValue $ws 'D22' '0.716'
$ws.Range('E22').Value = '  -7.77%  '
$ws.Range('E23').Value = '  -4.62%  '
Set-TextValue $ws 'D24' '2.48'
$ws.Range('E24').Value = '  +1.25%  '
Set-TextValue $ws 'D25' '13.00'
$ws.Range('E25').Value = '  -6.39%  '
Set-TextValue $ws 'D26' '81.38'
$ws.Range('E26').Value = '  -2.79%  '
Set-TextValue $ws 'D27' '10.61'
$ws.Range('E27').Value = '  +6.68%  '
$ws.Range('E28').Value = '  -0.31%  '
Set-TextValue $ws 'D29' '7.59'
$ws.Range('E29').Value = '  +0.93%  '
$ws.Range('E30').Value = '  -2.92%  '
$ws.Range('E31').Value = '  +0.00%  '
Set-TextValue $ws 'D32' '2.19'
$ws.Range('E32').Value = '  -3.97%  '
$ws.Range('E33').Value = '  -6.45%  '
Set-TextValue $ws 'D34' '27.30'
$ws.Range('E34').Value = '  -4.36%  '
$ws.Range('D35').Value = '0.0₃0844'
$ws.Range('E35').Value = '  -6.63%  '
Set-TextValue $ws 'D36' '1.06'
$ws.Range('E36').Value = '  -2.79%  '
$ws.Range('B37').Value = 'Filecoin'
$ws.Range('C37').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue $ws 'D37' '6.04'
$ws.Range('E37').Value = '  -4.98%  '
$ws.Range('B38').Value = 'dogwifhat'
$ws.Range('C38').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue $ws 'D38' '3.30'
$ws.Range('E38').Value = '  -8.87%  '
$ws.Range('E39').Value = '  -5.17%  '
Set-TextValue $ws 'D40' '50.97'
$ws.Range('E40').Value = '  -1.35%  '
Set-TextValue $ws 'D41' '9.19'
$ws.Range('E41').Value = '  -3.17%  '
Set-TextValue $ws 'D42' '446.08'
$ws.Range('E42').Value = '  -6.59%  '
Set-TextValue $ws 'D43' '0.293'
$ws.Range('E43').Value = '  -3.23%  '
Set-TextValue $ws 'D44' '41.10'
$ws.Range('E44').Value = '  +6.11%  '
Set-TextValue $ws 'D45' '0.0366'
$ws.Range('E45').Value = '  -4.64%  '
$ws.Range('E46').Value = '  -0.02%  '
$ws.Range('D47').Value = '2.831.46'
$ws.Range('E47').Value = '  -4.44%  '
Set-TextValue $ws 'D48' '130.81'
$ws.Range('E48').Value = '  -0.74%  '
Set-TextValue $ws 'D49' '25.91'
$ws.Range('E49').Value = '  +0.32%  '
Set-TextValue $ws 'D50' '1.00'
$ws.Range('E50').Value = '  +0.07%  '
$ws.Range('E51').Value = '  -3.81%  '
